$wb = $excel.ActiveWorkbook

# Rename sheets
$wb.Worksheets.Item("total mortality").Name = "mortality rates"
$wb.Worksheets.Item("mortality").Name = "causes of death"
$wb.Worksheets.Item("RRStunting").Name = "RR death by stunting"
$wb.Worksheets.Item("RRWasting").Name = "RR death by wasting"
$wb.Worksheets.Item("RRBreastfeeding").Name = "RR death by breastfeeding"
$wb.Worksheets.Item("RR Death by Birth Outcome").Name = "RR death by birth outcome"
$wb.Worksheets.Item("OR appropriateBF by interv").Name = "OR correctBF by interventn"

# Update the column title on "OR stunting for complements"
$ws = $wb.Worksheets.Item("OR stunting for complements")
$ws.Range("A1").Value = "Food security & education"

# Make "causes of death" the active sheet, with D47 selected
$wsActive = $wb.Worksheets.Item("causes of death")
$wsActive.Activate() | Out-Null
$wsActive.Range("D47").Select() | Out-Null
